# This script applies the weekly odds update for the Jogos_da_Semana_FlashScore_2024-11-16 workbook.
# Summary of changes:
#  - Updated several odds values in existing rows (2, 4, 7, 15, 17, 22)
#  - Inserted a new match row (Uruguay - Nacional vs Maldonado) at row 23,
#    which pushes the previous row 23 (Venezuela - Rayo Zuliano vs Caracas) down to row 24
#  - Updated the odds for the shifted Venezuela match (now row 24)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simple value cells in rows 2, 4, 7, 15, 17, 22
# Row 2
$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.48

# Row 4
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("AG4").Value = 1000

# Row 7
$ws.Range("G7").Value = 1.5
$ws.Range("H7").Value = 4.5
$ws.Range("I7").Value = 5.25
$ws.Range("J7").Value = 2.05
$ws.Range("L7").Value = 5.5
$ws.Range("AA7").Value = 11
$ws.Range("AD7").Value = 9
$ws.Range("AH7").Value = 19
$ws.Range("AJ7").Value = 17
$ws.Range("AM7").Value = 41
$ws.Range("AP7").Value = 15
$ws.Range("AQ7").Value = 21
$ws.Range("AW7").Value = 7.5
$ws.Range("AX7").Value = 29

# Row 15
$ws.Range("G15").Value = 2.8
$ws.Range("I15").Value = 2.7
$ws.Range("M15").Value = 1.1
$ws.Range("N15").Value = 7
$ws.Range("AA15").Value = 29
$ws.Range("AH15").Value = 6.5
$ws.Range("AK15").Value = 26
$ws.Range("AN15").Value = 4.75

# Row 17
$ws.Range("M17").Value = 1.1
$ws.Range("N17").Value = 7
$ws.Range("O17").Value = 1.5
$ws.Range("P17").Value = 2.5
$ws.Range("Q17").Value = 2.6
$ws.Range("R17").Value = 1.48
$ws.Range("S17").Value = 1.57
$ws.Range("T17").Value = 2.25
$ws.Range("W17").Value = 5
$ws.Range("AL17").Value = 51
$ws.Range("AN17").Value = 3.5
$ws.Range("AQ17").Value = 41
$ws.Range("AT17").Value = 2.25
$ws.Range("AU17").Value = 10

# Row 22
$ws.Range("M22").Value = 1.07
$ws.Range("N22").Value = 8.5
$ws.Range("Y22").Value = 9
$ws.Range("AC22").Value = 8.5
$ws.Range("AE22").Value = 17
$ws.Range("AF22").Value = 51
$ws.Range("AP22").Value = 23
$ws.Range("AZ22").Value = 101

# Insert a new row at position 23 (this shifts the existing row 23 down to row 24)
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 (Uruguay - Nacional vs Maldonado)
$ws.Range("A23").Value = "8OARDK0B"
$ws.Range("B23").Value = "16/11/2024"
$ws.Range("C23").Value = "19:30"
$ws.Range("D23").Value = "URUGUAY - PRIMERA DIVISION"
$ws.Range("E23").Value = "Nacional"
$ws.Range("F23").Value = "Maldonado"
$ws.Range("G23").Value = 1.14
$ws.Range("H23").Value = 7.5
$ws.Range("I23").Value = 17
$ws.Range("J23").Value = 1.53
$ws.Range("K23").Value = 2.75
$ws.Range("L23").Value = 15
$ws.Range("M23").Value = 1.03
$ws.Range("N23").Value = 17
$ws.Range("O23").Value = 1.18
$ws.Range("P23").Value = 4.5
$ws.Range("Q23").Value = 1.6
$ws.Range("R23").Value = 2.3
$ws.Range("S23").Value = 1.29
$ws.Range("T23").Value = 3.5
$ws.Range("U23").Value = 2.75
$ws.Range("V23").Value = 1.4
$ws.Range("W23").Value = 7
$ws.Range("X23").Value = 5.5
$ws.Range("Y23").Value = 11
$ws.Range("Z23").Value = 6
$ws.Range("AA23").Value = 13
$ws.Range("AB23").Value = 41
$ws.Range("AC23").Value = 13
$ws.Range("AD23").Value = 15
$ws.Range("AE23").Value = 41
$ws.Range("AF23").Value = 151
$ws.Range("AG23").Value = 201
$ws.Range("AH23").Value = 29
$ws.Range("AI23").Value = 67
$ws.Range("AJ23").Value = 41
$ws.Range("AK23").Value = 251
$ws.Range("AL23").Value = 126
$ws.Range("AM23").Value = 101
$ws.Range("AN23").Value = 3
$ws.Range("AO23").Value = 4.75
$ws.Range("AP23").Value = 21
$ws.Range("AQ23").Value = 11
$ws.Range("AR23").Value = 41
$ws.Range("AS23").Value = 201
$ws.Range("AT23").Value = 3.5
$ws.Range("AU23").Value = 13
$ws.Range("AV23").Value = 101
$ws.Range("AW23").Value = 15
$ws.Range("AX23").Value = 67
$ws.Range("AY23").Value = 67
$ws.Range("AZ23").Value = 351
$ws.Range("BA23").Value = 351
$ws.Range("BB23").Value = 501
$ws.Range("BC23").Value = 51
$ws.Range("BD23").Value = 51

# Update the shifted row 24 values (Venezuela - Rayo Zuliano vs Caracas) to the new odds
$ws.Range("A24").Value = "rZtf1obm"
$ws.Range("B24").Value = "16/11/2024"
$ws.Range("C24").Value = "16:30"
$ws.Range("D24").Value = "VENEZUELA - LIGA FUTVE"
$ws.Range("E24").Value = "Rayo Zuliano"
$ws.Range("F24").Value = "Caracas"
$ws.Range("G24").Value = 2.32
$ws.Range("H24").Value = 3.1
$ws.Range("I24").Value = 3
$ws.Range("J24").Value = 2.82
$ws.Range("K24").Value = 2.1
$ws.Range("L24").Value = 3.45
$ws.Range("M24").Value = 1.03
$ws.Range("N24").Value = 6.7
$ws.Range("O24").Value = 1.34
$ws.Range("P24").Value = 2.77
$ws.Range("Q24").Value = 1.98
$ws.Range("R24").Value = 1.65
$ws.Range("S24").Value = 1.39
$ws.Range("T24").Value = 2.57
$ws.Range("U24").Value = 1.75
$ws.Range("V24").Value = 1.85
$ws.Range("W24").Value = 7.5
$ws.Range("X24").Value = 11.25
$ws.Range("Y24").Value = 9
$ws.Range("Z24").Value = 24
$ws.Range("AA24").Value = 19.5
$ws.Range("AB24").Value = 30
$ws.Range("AC24").Value = 8.5
$ws.Range("AD24").Value = 6
$ws.Range("AE24").Value = 14
$ws.Range("AF24").Value = 70
$ws.Range("AG24").Value = 600
$ws.Range("AH24").Value = 8.5
$ws.Range("AI24").Value = 15
$ws.Range("AJ24").Value = 10.75
$ws.Range("AK24").Value = 37
$ws.Range("AL24").Value = 27
$ws.Range("AM24").Value = 37
$ws.Range("AN24").Value = 4.25
$ws.Range("AO24").Value = 11.5
$ws.Range("AP24").Value = 18
$ws.Range("AQ24").Value = 45
$ws.Range("AR24").Value = 70
$ws.Range("AS24").Value = 200
$ws.Range("AT24").Value = 2.6
$ws.Range("AU24").Value = 6.6
$ws.Range("AV24").Value = 55
$ws.Range("AW24").Value = 4.9
$ws.Range("AX24").Value = 16
$ws.Range("AY24").Value = 22
$ws.Range("AZ24").Value = 75
$ws.Range("BA24").Value = 100
$ws.Range("BB24").Value = 250
$ws.Range("BC24").Value = 51
$ws.Range("BD24").Value = 51
